$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert two new columns before column D (shifts D:K -> F:M) ---
$ws.Range("D1:E1").EntireColumn.Insert()

# --- Step 2: Copy number formats from column F (the old column D, now shifted)
#             onto the new D:E columns so the new cells reuse the existing
#             date / number styles instead of picking up column C's style. ---
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Give the two new columns a sensible width matching the other data columns.
$dataWidth = $ws.Range("F1").ColumnWidth
$ws.Range("D1").ColumnWidth = $dataWidth
$ws.Range("E1").ColumnWidth = $dataWidth

# --- Step 3: Populate the new D/E columns with the latest two quarters ---
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 248700
$ws.Range("E8").Value = 233200
$ws.Range("D9").Value = 82800
$ws.Range("E9").Value = 78600
$ws.Range("D10").Value = 165900
$ws.Range("E10").Value = 154600
$ws.Range("D12").Value = 71000
$ws.Range("E12").Value = 71100
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 100
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 248000
$ws.Range("E17").Value = 228500
$ws.Range("D18").Value = 700
$ws.Range("E18").Value = 4700
$ws.Range("D20").Value = 4600
$ws.Range("E20").Value = 5500
$ws.Range("D21").Value = 18100
$ws.Range("E21").Value = 22700
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 5400
$ws.Range("E23").Value = 10200
$ws.Range("D24").Value = 4800
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 600
$ws.Range("E26").Value = 10200
$ws.Range("D27").Value = 600
$ws.Range("E27").Value = 10200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -4600
$ws.Range("E32").Value = -5500
$ws.Range("D33").Value = 600
$ws.Range("E33").Value = 10200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 600
$ws.Range("E35").Value = 10200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 545000
$ws.Range("E41").Value = 236400
$ws.Range("D42").Value = 36200
$ws.Range("E42").Value = 183800
$ws.Range("D43").Value = 91600
$ws.Range("E43").Value = 107100
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 74400
$ws.Range("E45").Value = 43300
$ws.Range("D46").Value = 747300
$ws.Range("E46").Value = 570700
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 266600
$ws.Range("E48").Value = 266100
$ws.Range("D49").Value = 1052800
$ws.Range("E49").Value = 1075600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 80100
$ws.Range("E52").Value = 107100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2146700
$ws.Range("E54").Value = 2019400
$ws.Range("D57").Value = 84500
$ws.Range("E57").Value = 41900
$ws.Range("D58").Value = 100000
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 295400
$ws.Range("E59").Value = 261000
$ws.Range("D60").Value = 479800
$ws.Range("E60").Value = 302900
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 70300
$ws.Range("E62").Value = 100100
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 550100
$ws.Range("E66").Value = 403000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -1789700
$ws.Range("E72").Value = -1766000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1596600
$ws.Range("E76").Value = 1616400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 600
$ws.Range("E81").Value = 10200
$ws.Range("D83").Value = 12700
$ws.Range("E83").Value = 12500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 89900
$ws.Range("E89").Value = 41100
$ws.Range("D91").Value = -4000
$ws.Range("E91").Value = -3800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 144100
$ws.Range("E94").Value = -4600
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 75900
$ws.Range("E100").Value = -9000
$ws.Range("D101").Value = -1400
$ws.Range("E101").Value = -1200
$ws.Range("D102").Value = 308500
$ws.Range("E102").Value = 26300

# --- Step 4: A few historical cells were corrected in the same edit
#             (not just shifted) -- pin down the exact final values for
#             those rows across the whole D:M range. ---
$ws.Range("D58").Value = 100000
$ws.Range("E58").Value = "NA"
$ws.Range("F58").Value = "NA"
$ws.Range("G58").Value = "NA"
$ws.Range("H58").Value = "NA"
$ws.Range("I58").Value = "NA"
$ws.Range("J58").Value = "NA"
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 0
$ws.Range("D91").Value = -4000
$ws.Range("E91").Value = -3800
$ws.Range("F91").Value = -2300
$ws.Range("G91").Value = -1400
$ws.Range("H91").Value = -3100
$ws.Range("I91").Value = -2700
$ws.Range("J91").Value = -1900
$ws.Range("K91").Value = -4600
$ws.Range("L91").Value = -3700
$ws.Range("M91").Value = -2700
$ws.Range("D94").Value = 144100
$ws.Range("E94").Value = -4600
$ws.Range("F94").Value = -131500
$ws.Range("G94").Value = 11000
$ws.Range("H94").Value = -160900
$ws.Range("I94").Value = -260500
$ws.Range("J94").Value = -1600
$ws.Range("K94").Value = -37200
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = 14900

Write-Output "edit complete"
